$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") '61.041.84'
Set-TextValue $ws.Range("E2") '  -2.16%  '
Set-TextValue $ws.Range("D3") '3.376.78'
Set-TextValue $ws.Range("E3") '  +0.32%  '
Set-TextValue $ws.Range("E4") '  +0.01%  '
Set-TextValue $ws.Range("D5") '572.32'
Set-TextValue $ws.Range("E5") '  +1.20%  '
Set-TextValue $ws.Range("D6") '135.91'
Set-TextValue $ws.Range("E6") '  +9.77%  '
Set-TextValue $ws.Range("D8") '3.376.66'
Set-TextValue $ws.Range("E8") '  +0.20%  '
Set-TextValue $ws.Range("E9") '  +1.55%  '
Set-TextValue $ws.Range("D10") '7.59'
Set-TextValue $ws.Range("E10") '  +5.25%  '
Set-TextValue $ws.Range("E11") '  +3.85%  '
Set-TextValue $ws.Range("E12") '  +4.19%  '
Set-TextValue $ws.Range("D13") '3.958.60'
Set-TextValue $ws.Range("E13") '  +0.60%  '
Set-TextValue $ws.Range("E14") '  +1.49%  '
Set-TextValue $ws.Range("E15") '  +2.62%  '
Set-TextValue $ws.Range("D16") '3.387.54'
Set-TextValue $ws.Range("E16") '  +0.75%  '
Set-TextValue $ws.Range("D17") '25.14'
Set-TextValue $ws.Range("E17") '  +3.76%  '
Set-TextValue $ws.Range("D18") '61.257.64'
Set-TextValue $ws.Range("E18") '  -1.96%  '
Set-TextValue $ws.Range("D19") '14.01'
Set-TextValue $ws.Range("E19") '  +7.81%  '
Set-TextValue $ws.Range("E20") '  +3.87%  '
Set-TextValue $ws.Range("D21") '9.41'
Set-TextValue $ws.Range("E21") '  +1.82%  '
Set-TextValue $ws.Range("D22") '373.86'
Set-TextValue $ws.Range("E22") '  +1.51%  '
Set-TextValue $ws.Range("D23") '0.567'
Set-TextValue $ws.Range("E23") '  +3.21%  '
Set-TextValue $ws.Range("D24") '3.514.59'
Set-TextValue $ws.Range("E24") '  +0.46%  '
Set-TextValue $ws.Range("E25") '  +0.08%  '
Set-TextValue $ws.Range("D26") '70.59'
Set-TextValue $ws.Range("E26") '  -0.43%  '
Set-TextValue $ws.Range("D27") '0.0000117'
Set-TextValue $ws.Range("E27") '  +12.68%  '
Set-TextValue $ws.Range("D28") '1.66'
Set-TextValue $ws.Range("E28") '  +22.98%  '
Set-TextValue $ws.Range("D29") '7.69'
Set-TextValue $ws.Range("E29") '  +13.36%  '
Set-TextValue $ws.Range("E30") '  +0.09%  '
Set-TextValue $ws.Range("D31") '8.11'
Set-TextValue $ws.Range("E31") '  +5.81%  '
Set-TextValue $ws.Range("D32") '2.15'
Set-TextValue $ws.Range("E32") '  +1.98%  '
Set-TextValue $ws.Range("E33") '  +4.71%  '
Set-TextValue $ws.Range("E34") '  -0.05%  '
Set-TextValue $ws.Range("D35") '3.410.04'
Set-TextValue $ws.Range("E35") '  +0.49%  '
Set-TextValue $ws.Range("D36") '23.40'
Set-TextValue $ws.Range("E36") '  +4.18%  '
Set-TextValue $ws.Range("E37") '  +10.04%  '
Set-TextValue $ws.Range("E38") '  +7.51%  '
Set-TextValue $ws.Range("D39") '6.92'
Set-TextValue $ws.Range("E39") '  +5.45%  '
Set-TextValue $ws.Range("D40") '162.89'
Set-TextValue $ws.Range("E40") '  -0.79%  '
Set-TextValue $ws.Range("D41") '0.0788'
Set-TextValue $ws.Range("E41") '  +6.32%  '
Set-TextValue $ws.Range("E42") '  +0.05%  '
Set-TextValue $ws.Range("D43") '4.41'
Set-TextValue $ws.Range("E43") '  +5.48%  '
Set-TextValue $ws.Range("E44") '  +14.72%  '
Set-TextValue $ws.Range("D45") '0.760'
Set-TextValue $ws.Range("E45") '  -0.19%  '
Set-TextValue $ws.Range("D46") '41.34'
Set-TextValue $ws.Range("E46") '  +0.73%  '
Set-TextValue $ws.Range("D47") '1.60'
Set-TextValue $ws.Range("E47") '  +5.75%  '
Set-TextValue $ws.Range("D48") '23.34'
Set-TextValue $ws.Range("E48") '  +4.12%  '
Set-TextValue $ws.Range("E49") '  +6.22%  '
Set-TextValue $ws.Range("D50") '22.99'
Set-TextValue $ws.Range("E50") '  +15.78%  '
Set-TextValue $ws.Range("D51") '0.896'
Set-TextValue $ws.Range("E51") '  +7.26%  '
